# Update countries & provincias Spain
# Refresh the COVID-19 "paises" dataset snapshot:
#  - bump the "Datos actualizados..." timestamp
#  - update several countries' case/death/recovery counters
#  - Colombia overtakes Malasia in total cases, so the two swap places
#    in the ranking (row 50 <-> row 51)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 01:22"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1034115
$ws.Range("C4").Value = 23759
$ws.Range("D4").Value = 140501
$ws.Range("E4").Value = 834502
$ws.Range("G4").Value = 2315
$ws.Range("H4").Value = 59112

# Brasil (row 14)
$ws.Range("B14").Value = 72899
$ws.Range("C14").Value = 6398
$ws.Range("E14").Value = 35292
$ws.Range("G14").Value = 520
$ws.Range("H14").Value = 5063

# Canada (row 15)
$ws.Range("B15").Value = 50026
$ws.Range("C15").Value = 1526
$ws.Range("D15").Value = 19190
$ws.Range("E15").Value = 27977
$ws.Range("G15").Value = 152
$ws.Range("H15").Value = 2859

# Chequia (row 45)
$ws.Range("B45").Value = 7504
$ws.Range("C45").Value = 59
$ws.Range("D45").Value = 2948
$ws.Range("E45").Value = 4329
$ws.Range("G45").Value = 4
$ws.Range("H45").Value = 227

# Colombia overtakes Malasia: row 50 becomes Colombia (new, bigger totals),
# row 51 becomes Malasia (its figures are unchanged, only its rank moved).
$ws.Range("A50").Value = "Colombia"
$ws.Range("B50").Value = 5949
$ws.Range("C50").Value = 352
$ws.Range("D50").Value = 1268
$ws.Range("E50").Value = 4412
$ws.Range("F50").Value = 118
$ws.Range("G50").Value = 16
$ws.Range("H50").Value = 269

$ws.Range("A51").Value = "Malasia"
$ws.Range("B51").Value = 5851
$ws.Range("C51").Value = 31
$ws.Range("D51").Value = 4032
$ws.Range("E51").Value = 1719
$ws.Range("F51").Value = 36
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 100

# Principado de Andorra (row 96)
$ws.Range("D96").Value = 398
$ws.Range("E96").Value = 304
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 41

# Martinica (row 135)
$ws.Range("D135").Value = 83
$ws.Range("E135").Value = 78
$ws.Range("F135").Value = 5

# Guayana Francesa (row 143)
$ws.Range("B143").Value = 125
$ws.Range("C143").Value = 14
$ws.Range("D143").Value = 93
$ws.Range("E143").Value = 31

# Eritrea (row 173)
$ws.Range("D173").Value = 19
$ws.Range("E173").Value = 20
